# Reorders the beverage list (rows 3-23) to match the new query output from
# excel_writer.py, and normalizes previously-blank alcohol_content cells to
# the literal text "None" for rows with no alcohol (matches excel_input lists).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,(2, "Papillion", 25.36, 44.99, 13.5)
    ,(3, "Russian River Pliny the Elder", 16, 18.99, 8)
    ,(4, "Iced Tea", 18, 1.79, "None")
    ,(5, "Les Jamelles Syrah", 25.36, 10.99, 14.2)
    ,(6, "Goose Island IPA", 16, 9.99, 5.9)
    ,(7, "Mondavi Chardonnay", 25.36, 9.99, 12)
    ,(8, "Old Nation M-43", 16, 14.99, 6.8)
    ,(9, "Water", 16.9, 0.99, "None")
    ,(10, "Tomaiolo Pinot Grigio", 25.36, 9.99, 12)
    ,(11, "Murphys Irish Stout", 16, 10.99, 5)
    ,(12, "Relax Riesling", 25.36, 10.99, 11)
    ,(13, "Milk", 32, 2.99, "None")
    ,(14, "OVZ Zinfandel", 25.36, 13.99, 13.5)
    ,(15, "Iced Coffee", 16, 2.49, "None")
    ,(16, "Lange Pinot Noir", 25.36, 22.99, 13)
    ,(17, "Alchemist Heady Topper", 16, 17.99, 4.5)
    ,(18, "Smithwicks", 16, 12.99, 5)
    ,(19, "Alamos Malbec", 25.36, 12.99, 14)
    ,(20, "Trapiche Malbec", 25.36, 14.99, 13.8)
    ,(21, "Weldwerks Juicy Bits", 16, 13.99, 8)
    ,(22, "Matua Sauvignon Blanc", 25.36, 13.99, 12.5)
    ,(23, "Gatorade", 20, 2.49, "None")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
